$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 571
$wsExhibit.Range("F3").Value = 243
$wsExhibit.Range("F5").Value = 726
$wsExhibit.Range("F8").Value = 138
$wsExhibit.Range("F9").Value = 237
$wsExhibit.Range("F11").Value = 5799
$wsExhibit.Range("F12").Value = 44
$wsExhibit.Range("F13").Value = 33
$wsExhibit.Range("F16").Value = 539
$wsExhibit.Range("F17").Value = 339
$wsExhibit.Range("F22").Value = 89
$wsExhibit.Range("F23").Value = 302
$wsExhibit.Range("F24").Value = 1005
$wsExhibit.Range("F26").Value = 1755
$wsExhibit.Range("F27").Value = 452

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 47
$wsShow.Range("G4").Value = 880
$wsShow.Range("F5").Value = 261
$wsShow.Range("F6").Value = 294

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 197

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 197
$wsAll.Range("F3").Value = 571
$wsAll.Range("F4").Value = 243
$wsAll.Range("F6").Value = 726
$wsAll.Range("F10").Value = 138
$wsAll.Range("F11").Value = 237
$wsAll.Range("F13").Value = 5799
$wsAll.Range("F14").Value = 44
$wsAll.Range("F15").Value = 33
$wsAll.Range("F19").Value = 539
$wsAll.Range("F20").Value = 339
$wsAll.Range("F22").Value = 47
$wsAll.Range("G22").Value = 880
$wsAll.Range("F25").Value = 261
$wsAll.Range("F26").Value = 294
$wsAll.Range("F32").Value = 89
$wsAll.Range("F33").Value = 302
$wsAll.Range("F34").Value = 1005
$wsAll.Range("F36").Value = 1755
$wsAll.Range("F37").Value = 452
